$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

$ws.Cells.Item($row, 1).Value = 46025
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat

$ws.Cells.Item($row, 2).Value = 169
$ws.Cells.Item($row, 3).Value = 179
$ws.Cells.Item($row, 4).Value = 171
